$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R32 moves from the "150" value group (row 25) to the "470" value group (row 24).
# Set row 25's designator list first (drop R32), then row 24's (add R32) so the
# new shared-string entries get appended to the sst table in the same order as
# the target workbook.
$ws.Range("A25").Value = "R3, R4, R5, R6, R15, R16, R24, R25, R35, R45, R80, R81, R82, R83, R84"
$ws.Range("A24").Value = "R1, R2, R7, R9, R12, R32, R36, R56, R57, R58, R59, R60, R61, R62, R63, R64, R65, R66, R67, R68, R69, R70, R71, R72, R73, R74, R75, R76, R77, R78, R79"

# Update the quantity counts to reflect the moved designator.
$ws.Range("D24").Value = 31
$ws.Range("D25").Value = 15

# Move the active selection/cursor to A33, matching the saved view state.
$ws.Range("A33").Select()
